# "fix bug test date": the hard-coded specimen-id suffix on the report
# header ("... 180890 -370843") is replaced with a template placeholder
# so the real value can be substituted at generation time.
$d = $word.ActiveDocument

$d.Content.Find.Execute("-370843", $true, $false, $false, $false, $false,
                         $true, 1, $false, "–{specimenId}", 2)
